$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-03-26"

# Update the row label for March to reflect the new "through" date
$ws.Range("A4").Value = "March (through 03-26)"

# Update March row (row 4) values
$ws.Range("C4").Value = 37
$ws.Range("D4").Value = 47
$ws.Range("E4").Value = 52
$ws.Range("G4").Value = 52
$ws.Range("H4").Value = 71
$ws.Range("I4").Value = 104

# Update Total row (row 5) values
$ws.Range("C5").Value = 124
$ws.Range("D5").Value = 178
$ws.Range("E5").Value = 189
$ws.Range("G5").Value = 193
$ws.Range("H5").Value = 413
$ws.Range("I5").Value = 404
